$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# YDS sheet: append Week 13 play-by-play yardage logs (run/pass, for/against)
# ------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

# B2 = OFF Run yards log
$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Text + " 2 1 3 17 9 1 0 1 15 2 2 3 2 3 8 17"

# C2 = DEF Run yards log (yards allowed)
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Text + " 0 5 9 2 2 1 4 9 7 18 9 4 1 3 0 7 6 2 2 6 0 0 1 1 3 -1 3 5 -1 8 34 1 4 2 2 7 12 3 2 0"

# B3 = OFF Pass yards log
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Text + " 8 7 3 10 7 5 10 29 4 19 5 18 14 19 11 -2 15 4 14 7 11"

# C3 = DEF Pass yards log (yards allowed)
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Text + " 16 9 36 12 10 12 5 25 13 16 12 22 -2 2 28 9 7 7 2 1"

# ------------------------------------------------------------------
# OFF sheet: update running totals for Home (row2) and Road (row3)
# ------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = [int]$wsOFF.Range("C2").Text + 7
$wsOFF.Range("D2").Value = [int]$wsOFF.Range("D2").Text + 1
$wsOFF.Range("F2").Value = [int]$wsOFF.Range("F2").Text + 1
$wsOFF.Range("G2").Value = [int]$wsOFF.Range("G2").Text + 7
$wsOFF.Range("L2").Value = [int]$wsOFF.Range("L2").Text + 37
$wsOFF.Range("M2").Value = [int]$wsOFF.Range("M2").Text + 23
$wsOFF.Range("O2").Value = [int]$wsOFF.Range("O2").Text + 3
$wsOFF.Range("P2").Value = [int]$wsOFF.Range("P2").Text + 2
$wsOFF.Range("Q2").Value = [int]$wsOFF.Range("Q2").Text + 57

$wsOFF.Range("C3").Value = [int]$wsOFF.Range("C3").Text + 16
$wsOFF.Range("E3").Value = [int]$wsOFF.Range("E3").Text + 2
$wsOFF.Range("F3").Value = [int]$wsOFF.Range("F3").Text + 9
$wsOFF.Range("I3").Value = [int]$wsOFF.Range("I3").Text + 5
$wsOFF.Range("J3").Value = [int]$wsOFF.Range("J3").Text + 4
$wsOFF.Range("N3").Value = [int]$wsOFF.Range("N3").Text + 2

# ------------------------------------------------------------------
# DEF sheet: update running totals for Home (row2) and Road (row3)
# ------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value = [int]$wsDEF.Range("C2").Text + 23
$wsDEF.Range("D2").Value = [int]$wsDEF.Range("D2").Text + 1
$wsDEF.Range("E2").Value = [int]$wsDEF.Range("E2").Text + 1
$wsDEF.Range("F2").Value = [int]$wsDEF.Range("F2").Text + 6
$wsDEF.Range("G2").Value = [int]$wsDEF.Range("G2").Text + 5
$wsDEF.Range("I2").Value = [int]$wsDEF.Range("I2").Text + 1
$wsDEF.Range("J2").Value = [int]$wsDEF.Range("J2").Text + 3
$wsDEF.Range("L2").Value = [int]$wsDEF.Range("L2").Text + 25
$wsDEF.Range("M2").Value = [int]$wsDEF.Range("M2").Text + 20
$wsDEF.Range("O2").Value = [int]$wsDEF.Range("O2").Text + 1
$wsDEF.Range("P2").Value = [int]$wsDEF.Range("P2").Text + 1
$wsDEF.Range("Q2").Value = [int]$wsDEF.Range("Q2").Text + 72

$wsDEF.Range("C3").Value = [int]$wsDEF.Range("C3").Text + 6
$wsDEF.Range("E3").Value = [int]$wsDEF.Range("E3").Text + 1
$wsDEF.Range("F3").Value = [int]$wsDEF.Range("F3").Text + 9
$wsDEF.Range("G3").Value = [int]$wsDEF.Range("G3").Text + 2
$wsDEF.Range("H3").Value = [int]$wsDEF.Range("H3").Text + 2
$wsDEF.Range("I3").Value = [int]$wsDEF.Range("I3").Text + 2
$wsDEF.Range("J3").Value = [int]$wsDEF.Range("J3").Text + 5
$wsDEF.Range("N3").Value = [int]$wsDEF.Range("N3").Text + 2

# ------------------------------------------------------------------
# ST sheet: special teams totals + append kick/punt distance logs
# ------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = [int]$wsST.Range("B2").Text + 4
$wsST.Range("D2").Value = [int]$wsST.Range("D2").Text + 1
$wsST.Range("H2").Value = [int]$wsST.Range("H2").Text + 1
$wsST.Range("B3").Value = [int]$wsST.Range("B3").Text + 1

# B4 = KO distance log
$wsST.Range("B4").Value = $wsST.Range("B4").Text + " 55 64 62"
# B5 = KO return log
$wsST.Range("B5").Value = $wsST.Range("B5").Text + " 21 5 13"
# B6 = PT distance log
$wsST.Range("B6").Value = $wsST.Range("B6").Text + " 79 28 32"
# D3 = opponent KO distance log
$wsST.Range("D3").Value = $wsST.Range("D3").Text + " 33"
# D4 = opponent KO return log
$wsST.Range("D4").Value = $wsST.Range("D4").Text + " 0"
# D5 = opponent PT distance log
$wsST.Range("D5").Value = $wsST.Range("D5").Text + " 0"

# ------------------------------------------------------------------
# TURNS sheet: turnover totals
# ------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("B2").Value = [int]$wsTURNS.Range("B2").Text + 1
$wsTURNS.Range("E2").Value = [int]$wsTURNS.Range("E2").Text + 2
$wsTURNS.Range("D3").Value = [int]$wsTURNS.Range("D3").Text - 1
$wsTURNS.Range("E3").Value = [int]$wsTURNS.Range("E3").Text - 1

# ------------------------------------------------------------------
# PEN sheet: penalty totals
# ------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("D2").Value = [int]$wsPEN.Range("D2").Text + 2
$wsPEN.Range("D4").Value = [int]$wsPEN.Range("D4").Text + 1
